# Update the crypto price/volume table with the latest scraped figures.
# Rows 30 and 31 (Fetch.AI / PancakeSwap) also swap rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.537.65'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '3.508.83'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'609.47"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = "'152.38"
$ws.Range("E6").Value = '  +1.04%  '
$ws.Range("D7").Value = '3.506.76'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("E11").Value = '  +8.48%  '
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").Value = "'32.81"
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").Value = "'0.0000217"
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '4.098.15'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '3.507.37'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '67.395.34'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = "'6.59"
$ws.Range("E19").Value = '  +2.87%  '
$ws.Range("D20").Value = "'15.61"
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("E21").Value = '  +7.01%  '
$ws.Range("D22").Value = "'448.89"
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("D24").Value = "'78.29"
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = '3.647.13'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = "'8.89"
$ws.Range("E28").Value = '  +6.98%  '
$ws.Range("D29").Value = "'10.15"
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.52"
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = "'1.66"
$ws.Range("E31").Value = '  +7.42%  '
$ws.Range("D32").Value = "'0.169"
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = "'25.78"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("D37").Value = '3.501.11'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = "'0.0900"
$ws.Range("E42").Value = '  +3.03%  '
$ws.Range("D43").Value = "'173.47"
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("D44").Value = "'5.51"
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").Value = "'30.20"
$ws.Range("E45").Value = '  +10.71%  '
$ws.Range("D46").Value = "'0.884"
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("E47").Value = '  +3.04%  '
$ws.Range("D48").Value = "'1.31"
$ws.Range("E48").Value = '  +4.12%  '
$ws.Range("D49").Value = "'7.68"
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("D50").Value = "'2.52"
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("D51").Value = "'0.255"
$ws.Range("E51").Value = '  +3.67%  '
